# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")
$ws.Activate()

# Data refresh: regcntr_id values changed in the master-reg_center_user sheet
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Restore the sheet's scroll position / selection as left by the author after
# the refresh (best-effort — not all view-state properties persist).
$ws.Range("C19").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
